$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 445/446, pushing the existing rows 445-517 down to 447-519.
$ws.Rows("445:446").Insert()

# --- New row 445 ---
$ws.Cells.Item(445, 1).Value = 7
$ws.Cells.Item(445, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(445, 3).Value = "Ñuble"
$ws.Cells.Item(445, 4).Value = 45180
$ws.Cells.Item(445, 5).Value = 16
$ws.Cells.Item(445, 6).Value = 100114013
$ws.Cells.Item(445, 7).Value = "Zanahoria"
$ws.Cells.Item(445, 8).Value = "Sin especificar"
$ws.Cells.Item(445, 9).Value = "Primera"
$ws.Cells.Item(445, 10).Value = 120
$ws.Cells.Item(445, 11).Value = 6000
$ws.Cells.Item(445, 12).Value = 6000
$ws.Cells.Item(445, 13).Value = 6000
$ws.Cells.Item(445, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(445, 15).Value = "Región de Ñuble"
$ws.Cells.Item(445, 16).Value = 300
$ws.Cells.Item(445, 17).Value = 20
$ws.Cells.Item(445, 18).Value = "Hortaliza"

# --- New row 446 ---
$ws.Cells.Item(446, 1).Value = 7
$ws.Cells.Item(446, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(446, 3).Value = "Ñuble"
$ws.Cells.Item(446, 4).Value = 45180
$ws.Cells.Item(446, 5).Value = 16
$ws.Cells.Item(446, 6).Value = 100114013
$ws.Cells.Item(446, 7).Value = "Zanahoria"
$ws.Cells.Item(446, 8).Value = "Sin especificar"
$ws.Cells.Item(446, 9).Value = "Segunda"
$ws.Cells.Item(446, 10).Value = 80
$ws.Cells.Item(446, 11).Value = 5000
$ws.Cells.Item(446, 12).Value = 5000
$ws.Cells.Item(446, 13).Value = 5000
$ws.Cells.Item(446, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(446, 15).Value = "Región de Ñuble"
$ws.Cells.Item(446, 16).Value = 250
$ws.Cells.Item(446, 17).Value = 20
$ws.Cells.Item(446, 18).Value = "Hortaliza"

Write-Host ("UsedRange: " + $ws.UsedRange.Address())
